$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F2:F11 "想去人数" counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 195
$ws1.Range("F3").Value = 521
$ws1.Range("F4").Value = 36
$ws1.Range("F5").Value = 24
$ws1.Range("F6").Value = 15
$ws1.Range("F7").Value = 33
$ws1.Range("F8").Value = 26
$ws1.Range("F9").Value = 249
$ws1.Range("F10").Value = 2758
$ws1.Range("F11").Value = 28

# Sheet "全部类型" (sheet4): update same counts (row 3 unaffected)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 195
$ws4.Range("F4").Value = 521
$ws4.Range("F5").Value = 36
$ws4.Range("F6").Value = 24
$ws4.Range("F7").Value = 15
$ws4.Range("F8").Value = 33
$ws4.Range("F9").Value = 26
$ws4.Range("F10").Value = 249
$ws4.Range("F11").Value = 2758
$ws4.Range("F12").Value = 28
